# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" / price figures to both the "展览" and
# "全部类型" worksheets (which carry duplicated data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("G2").Value2  = 65
    $ws.Range("F3").Value2  = 6389
    $ws.Range("F5").Value2  = 22
    $ws.Range("F7").Value2  = 1918
    $ws.Range("F8").Value2  = 1452
    $ws.Range("F10").Value2 = 978
    $ws.Range("F11").Value2 = 290
    $ws.Range("F12").Value2 = 5595
    $ws.Range("F13").Value2 = 73
}
